# Adds a new weekly price report (fecha 44468) for
# "Comercializadora del Agro de Limarí" / Pepino dulce, inserted right
# before the existing 44308 block (row 212), pushing all subsequent
# rows down by 4 and growing the used range from R264 to R268.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 212 (shifts old rows 212:264 down to 216:268,
# inheriting formatting from the row that was at 212, same as Excel's
# normal Insert behaviour).
$ws.Rows("212:215").Insert()

# Quality tiers for the new 44468 report, in on-sheet order.
$newRows = @(
    @(212, "Especial", 240, 13500, 14000, 13750, 764),
    @(213, "Primera",  400, 11500, 12000, 11750, 653),
    @(214, "Segunda",  300,  9500, 10000,  9750, 542),
    @(215, "Tercera",  200,  6500,  7000,  6750, 375)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Range("A$row").Value = 2
    $ws.Range("B$row").Value = "Comercializadora del Agro de Limarí"
    $ws.Range("C$row").Value = "Coquimbo"
    $ws.Range("D$row").Value = 44468
    $ws.Range("E$row").Value = 4
    $ws.Range("F$row").Value = 100112043
    $ws.Range("G$row").Value = "Pepino dulce"
    $ws.Range("H$row").Value = "Cultivar IV Región"
    $ws.Range("I$row").Value = $r[1]
    $ws.Range("J$row").Value = $r[2]
    $ws.Range("K$row").Value = $r[3]
    $ws.Range("L$row").Value = $r[4]
    $ws.Range("M$row").Value = $r[5]
    $ws.Range("N$row").Value = "$/bandeja 18 kilos"
    $ws.Range("O$row").Value = "Provincia de Limarí"
    $ws.Range("P$row").Value = $r[6]
    $ws.Range("Q$row").Value = 18
    $ws.Range("R$row").Value = "Hortaliza"
}
